# Apply "MeritOrderLänder" simplification updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column F (Steinkohle_0): 4446.0 -> 3276.0 for rows 2-11
$ws.Range("F2:F11").Value = 3276.0

# Column N (Wasserkraft): 950.0 -> 1000.0 for rows 2-11
$ws.Range("N2:N11").Value = 1000.0

# Column I (Erdgas_0): new values per row (rows 2-11)
$newI = @(5039.1, 4923.1, 4573.1, 4214.1, 3904.1, 3854.1, 4068.1, 4310.1, 4645.1, 5179.1)
for ($i = 0; $i -lt $newI.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $newI[$i]
}
